# Applies the "Automatic update of files" edit to the Artfynd sheet:
# the observation records in rows 5-9 get their data fields re-shuffled
# (row5<->row6 swap; row7->row9->row8->row7 cycle) while a handful of
# location/observer columns that are identical across these rows stay put.
# Below are the concrete per-cell before/after values needed to reach the
# target state, applied directly so no row-move ambiguity remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (becomes the old row 6 record: Harticka / Pelloporus leporinus) ---
$ws.Range("A5").Value = 112038436
$ws.Range("B5").Value = 89401
$ws.Range("E5").Value = 1108
$ws.Range("F5").Value = "Harticka"
$ws.Range("G5").Value = "Pelloporus leporinus"
$ws.Range("H5").Value = "(Fr.) Krieglst."
$ws.Range("Q5").Value = 515951.3091604927
$ws.Range("R5").Value = 7184319.58691278
$ws.Range("Z5").Value = "13:28"
$ws.Range("AB5").Value = "13:28"
$ws.Range("AJ5").Value = "gran"
$ws.Range("AK5").Value = "Picea abies"
$ws.Range("AM5").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO5").Value = "Standing dead tree/snags # Picea abies"

# --- Row 6 (becomes the old row 5 record: Garnlav / Alectoria sarmentosa) ---
$ws.Range("A6").Value = 112035549
$ws.Range("B6").Value = 77515
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("Q6").Value = 515977.3292799139
$ws.Range("R6").Value = 7184566.677681392
$ws.Range("Z6").Value = "10:51"
$ws.Range("AB6").Value = "10:51"
$ws.Range("AJ6").Value = ""
$ws.Range("AK6").Value = ""
$ws.Range("AM6").Value = "Gren på levande träd"
$ws.Range("AO6").Value = "Branch on living tree"

# --- Row 7 (becomes the old row 8 record: Rosenticka / Rhodofomes roseus) ---
$ws.Range("A7").Value = 112038473
$ws.Range("B7").Value = 89686
$ws.Range("E7").Value = 658
$ws.Range("F7").Value = "Rosenticka"
$ws.Range("G7").Value = "Rhodofomes roseus"
$ws.Range("H7").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
# "4" must land as text (matches source inlineStr), not get auto-coerced to a number
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "4"
$ws.Range("I7").ClearFormats()
$ws.Range("J7").Value = "fruktkroppar"
$ws.Range("Q7").Value = 516057.2181607572
$ws.Range("R7").Value = 7184319.723381012
$ws.Range("Z7").Value = "13:34"
$ws.Range("AB7").Value = "13:34"

# --- Row 8 (becomes the old row 9 record: Garnlav / Alectoria sarmentosa) ---
$ws.Range("A8").Value = 112038529
$ws.Range("B8").Value = 77515
$ws.Range("E8").Value = 6425
$ws.Range("F8").Value = "Garnlav"
$ws.Range("G8").Value = "Alectoria sarmentosa"
$ws.Range("H8").Value = "(Ach.) Ach."
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = ""
$ws.Range("Q8").Value = 515871.5299412137
$ws.Range("R8").Value = 7184628.386151251
$ws.Range("Z8").Value = "14:07"
$ws.Range("AB8").Value = "14:07"
$ws.Range("AH8").Value = "Gransumpskog"
$ws.Range("AJ8").Value = ""
$ws.Range("AK8").Value = ""
$ws.Range("AM8").Value = ""
$ws.Range("AO8").Value = ""

# --- Row 9 (becomes the old row 7 record: Ullticka / Phellinidium ferrugineofuscum) ---
$ws.Range("A9").Value = 112038134
$ws.Range("B9").Value = 89405
$ws.Range("E9").Value = 1202
$ws.Range("F9").Value = "Ullticka"
$ws.Range("G9").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H9").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q9").Value = 515925.2595200292
$ws.Range("R9").Value = 7184319.449006356
$ws.Range("Z9").Value = "13:27"
$ws.Range("AB9").Value = "13:27"
$ws.Range("AH9").Value = "Blåbärsgranskog"
$ws.Range("AJ9").Value = "gran"
$ws.Range("AK9").Value = "Picea abies"
$ws.Range("AM9").Value = "Liggande död trädstam, markontakt"
$ws.Range("AO9").Value = "Horizontal, dead with ground contact # Picea abies"
